$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E (shifts Engine No. etc. one to the right)
$ws.Range("E1").EntireColumn.Insert()

# Set header text for the newly inserted column E
$ws.Range("E1").Value = "PUC Amount"

# Set the width for the new column E
$ws.Columns.Item(5).ColumnWidth = 18.284615384615385

# Update the selection to match the target state
$ws.Range("E1:E1048576").Select()
